$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.429598666666667
$ws.Range("H2").Value = 22.288796
$ws.Range("I2").Value = 0.2633764192298049
$ws.Range("J2").Value = 0.2633764192298049
$ws.Range("M2").Value = 6.663840333333333
$ws.Range("N2").Value = 19.991521
$ws.Range("O2").Value = 0.3746160267057107
$ws.Range("P2").Value = 0.3746160267057107
$ws.Range("Q2").Value = 49.50965925541289
$ws.Range("R2").Value = 445.586933298716
$ws.Range("S2").Value = 0.09866502769984704
$ws.Range("T2").Value = 0.09866502769984704
$ws.Range("G3").Value = 7.429598666666667
$ws.Range("H3").Value = 22.288796
$ws.Range("I3").Value = 0.2633764192298049
$ws.Range("J3").Value = 0.2633764192298049
$ws.Range("M3").Value = 11.12461466666667
$ws.Range("N3").Value = 33.373844
$ws.Range("O3").Value = 0.6253839732942893
$ws.Range("P3").Value = 0.6253839732942893
$ws.Range("Q3").Value = 82.65142229464711
$ws.Range("R3").Value = 743.862800651824
$ws.Range("S3").Value = 0.1647113915299578
$ws.Range("T3").Value = 0.1647113915299578
$ws.Range("I4").Value = 0.07089061759860023
$ws.Range("J4").Value = 0.07089061759860024
$ws.Range("M4").Value = 6.663840333333333
$ws.Range("N4").Value = 19.991521
$ws.Range("O4").Value = 0.3746160267057107
$ws.Range("P4").Value = 0.3746160267057107
$ws.Range("Q4").Value = 13.32606135346567
$ws.Range("R4").Value = 119.934552181191
$ws.Range("S4").Value = 0.02655676149550155
$ws.Range("T4").Value = 0.02655676149550155
$ws.Range("I5").Value = 0.07089061759860023
$ws.Range("J5").Value = 0.07089061759860024
$ws.Range("M5").Value = 11.12461466666667
$ws.Range("N5").Value = 33.373844
$ws.Range("O5").Value = 0.6253839732942893
$ws.Range("P5").Value = 0.6253839732942893
$ws.Range("Q5").Value = 22.24652605196933
$ws.Range("R5").Value = 200.218734467724
$ws.Range("S5").Value = 0.04433385610309868
$ws.Range("T5").Value = 0.04433385610309869
$ws.Range("G6").Value = 7.761126333333333
$ws.Range("H6").Value = 23.283379
$ws.Range("I6").Value = 0.2751289476825233
$ws.Range("J6").Value = 0.2751289476825233
$ws.Range("M6").Value = 6.663840333333333
$ws.Range("N6").Value = 19.991521
$ws.Range("O6").Value = 0.3746160267057107
$ws.Range("P6").Value = 0.3746160267057107
$ws.Range("Q6").Value = 51.71890669216211
$ws.Range("R6").Value = 465.470160229459
$ws.Range("S6").Value = 0.1030677132125502
$ws.Range("T6").Value = 0.1030677132125502
$ws.Range("G7").Value = 7.761126333333333
$ws.Range("H7").Value = 23.283379
$ws.Range("I7").Value = 0.2751289476825233
$ws.Range("J7").Value = 0.2751289476825233
$ws.Range("M7").Value = 11.12461466666667
$ws.Range("N7").Value = 33.373844
$ws.Range("O7").Value = 0.6253839732942893
$ws.Range("P7").Value = 0.6253839732942893
$ws.Range("Q7").Value = 86.33953983765288
$ws.Range("R7").Value = 777.055858538876
$ws.Range("S7").Value = 0.172061234469973
$ws.Range("T7").Value = 0.1720612344699731
$ws.Range("G8").Value = 0.9721176666666667
$ws.Range("H8").Value = 2.916353
$ws.Range("I8").Value = 0.03446119791980235
$ws.Range("J8").Value = 0.03446119791980236
$ws.Range("M8").Value = 6.663840333333333
$ws.Range("N8").Value = 19.991521
$ws.Range("O8").Value = 0.3746160267057107
$ws.Range("P8").Value = 0.3746160267057107
$ws.Range("Q8").Value = 6.478036915879222
$ws.Range("R8").Value = 58.302332242913
$ws.Range("S8").Value = 0.01290971704023546
$ws.Range("T8").Value = 0.01290971704023546
$ws.Range("G9").Value = 0.9721176666666667
$ws.Range("H9").Value = 2.916353
$ws.Range("I9").Value = 0.03446119791980235
$ws.Range("J9").Value = 0.03446119791980236
$ws.Range("M9").Value = 11.12461466666667
$ws.Range("N9").Value = 33.373844
$ws.Range("O9").Value = 0.6253839732942893
$ws.Range("P9").Value = 0.6253839732942893
$ws.Range("Q9").Value = 10.81443445232578
$ws.Range("R9").Value = 97.32991007093199
$ws.Range("S9").Value = 0.02155148087956689
$ws.Range("T9").Value = 0.0215514808795669
$ws.Range("G10").Value = 10.04645066666667
$ws.Range("H10").Value = 30.139352
$ws.Range("I10").Value = 0.3561428175692692
$ws.Range("J10").Value = 0.3561428175692692
$ws.Range("M10").Value = 6.663840333333333
$ws.Range("N10").Value = 19.991521
$ws.Range("O10").Value = 0.3746160267057107
$ws.Range("P10").Value = 0.3746160267057107
$ws.Range("Q10").Value = 66.94794315937689
$ws.Range("R10").Value = 602.5314884343919
$ws.Range("S10").Value = 0.1334168072575764
$ws.Range("T10").Value = 0.1334168072575764
$ws.Range("G11").Value = 10.04645066666667
$ws.Range("H11").Value = 30.139352
$ws.Range("I11").Value = 0.3561428175692692
$ws.Range("J11").Value = 0.3561428175692692
$ws.Range("M11").Value = 11.12461466666667
$ws.Range("N11").Value = 33.373844
$ws.Range("O11").Value = 0.6253839732942893
$ws.Range("P11").Value = 0.6253839732942893
$ws.Range("Q11").Value = 111.7628924343431
$ws.Range("R11").Value = 1005.866031909088
$ws.Range("S11").Value = 0.2227260103116928
$ws.Range("T11").Value = 0.2227260103116928
